$d = $word.ActiveDocument

# --- "Common" policy table: PTYM-COM001 / COM002 / COM003 ---
# The table currently has 4 rows (header + 3 policy rows):
#   PTYM-COM001 | Read history on all PTYM that they can read
#   PTYM-COM002 | Read all the PTYM concerning the current party
#   PTYM-COM003 | Read all the history of PTYM concerning the current party
# After the edit it has only 2 policy rows, each having absorbed the text
# of the following one, and the old COM003 row is removed entirely:
#   PTYM-COM001 | Read all the PTYM concerning the current party
#   PTYM-COM002 | Read all the history of PTYM concerning the current party
$commonTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(2,1).Range.Text -match "PTYM-COM001") {
        $commonTable = $candidate
        break
    }
}

$commonTable.Cell(2,2).Range.Text = "Read all the PTYM concerning the current party"
$commonTable.Cell(3,2).Range.Text = "Read all the history of PTYM concerning the current party"
$commonTable.Rows.Item(4).Delete()

# --- "Organisation" policy table: add PTYM-ORG002 row ---
$orgTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Cell(2,1).Range.Text -match "PTYM-ORG001") {
        $orgTable = $candidate
        break
    }
}

$newRow = $orgTable.Rows.Add()
$newRowIndex = $newRow.Index
$orgTable.Cell($newRowIndex, 1).Range.Text = "PTYM-ORG002"
$orgTable.Cell($newRowIndex, 2).Range.Text = "Read PTYM history on all parties owned by the entity owning the organisation party."
$orgTable.Cell($newRowIndex, 3).Range.Text = "DONE"
